$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new "Preferred" columns (R = PreferredStartDate, S = PreferredEndDate) ---
# (values are written in the same order Excel originally appended them to the shared
# string table: the two new date strings, then the two new header labels, then the two
# updated reference-id values below)
$ws.Range("R2").Value = "2017-09-10 10:00:00"
$ws.Range("S2").Value = "2017-09-10 16:00:00"
$ws.Range("R1").Value = "PreferredStartDate"
$ws.Range("S1").Value = "PreferredEndDate"

# --- Update the two reference-id cells whose values changed ---
$ws.Range("D2").Value = "a0Nq0000003PF2Z"
$ws.Range("D3").Value = "a0Nq0000003PF2e"

# Match the formatting of the existing analogous columns: header cells (R1/S1) copy the
# plain wrap-text header format used by N1/O1; the data cells (R2/S2) copy the
# left-aligned text-number-format used by L2/M2.
$ws.Range("N1:O1").Copy()
$ws.Range("R1:S1").PasteSpecial(-4122)

$ws.Range("L2:M2").Copy()
$ws.Range("R2:S2").PasteSpecial(-4122)

# --- Column widths for the new columns ---
$ws.Columns.Item(18).ColumnWidth = 41.498697916666664
$ws.Columns.Item(19).ColumnWidth = 27.666666666666668

# --- Selection moved from N3 to D3, and the frozen/scrolled topLeftCell (L1) is cleared ---
$ws.Range("D3").Select()
